# Generate Report for Handoff
# 9fb53365-d646-45a2-8ef4-49cdaa28c62b.md is now ready for handoff and moves
# to row 2 (in sync status), while 5fe23c8a-ea68-4f02-8d8b-6baca2a7de74.md
# moves to row 3 and becomes "Ready for handoff" with an updated handback
# datetime and a handback-version-mismatch error detail.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "9fb53365-d646-45a2-8ef4-49cdaa28c62b.md"
$ov.Range("B2").Value = "e2e\9fb53365-d646-45a2-8ef4-49cdaa28c62b.md"

$ov.Range("A3").Value = "5fe23c8a-ea68-4f02-8d8b-6baca2a7de74.md"
$ov.Range("B3").Value = "e2e\5fe23c8a-ea68-4f02-8d8b-6baca2a7de74.md"
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-09-03 06:51:40"

# Keep the hyperlink display text (and cached display cache) in sync with
# the swapped file names, same way Excel does when you edit a linked cell
# via the Hyperlinks API.
$ov.Hyperlinks.Item(1).TextToDisplay = "e2e\9fb53365-d646-45a2-8ef4-49cdaa28c62b.md"
$ov.Hyperlinks.Item(2).TextToDisplay = "e2e\5fe23c8a-ea68-4f02-8d8b-6baca2a7de74.md"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "9fb53365-d646-45a2-8ef4-49cdaa28c62b.md"
$zh.Range("G2").Value = "9fb53365-d646-45a2-8ef4-49cdaa28c62b.fca86b0be0bc02ca2aa80548d8c5776cca403c63.zh-cn.xlf"
$zh.Range("I2").Value = "9fb53365-d646-45a2-8ef4-49cdaa28c62b.md"
$zh.Range("J2").Value = "9fb53365-d646-45a2-8ef4-49cdaa28c62b.fca86b0be0bc02ca2aa80548d8c5776cca403c63.zh-cn.xlf"

$zh.Range("A3").Value = "5fe23c8a-ea68-4f02-8d8b-6baca2a7de74.md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("G3").Value = "5fe23c8a-ea68-4f02-8d8b-6baca2a7de74.e0e312a260014abee75e5e6bf1ab9bc9e18a2d7d.zh-cn.xlf"
$zh.Range("H3").Value = "2016-09-03 06:51:35"
$zh.Range("I3").Value = "5fe23c8a-ea68-4f02-8d8b-6baca2a7de74.md"
$zh.Range("J3").Value = "5fe23c8a-ea68-4f02-8d8b-6baca2a7de74.e0e312a260014abee75e5e6bf1ab9bc9e18a2d7d.zh-cn.xlf"
$zh.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/470930ec7e4338bec95c5b6e9ad25607e2184c5a/e2e/5fe23c8a-ea68-4f02-8d8b-6baca2a7de74.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7442d1f7ad367e65fbeed7c5c61faebaeb5310d8/e2e/5fe23c8a-ea68-4f02-8d8b-6baca2a7de74.md."

$zh.Hyperlinks.Item(1).TextToDisplay = "9fb53365-d646-45a2-8ef4-49cdaa28c62b.md"
$zh.Hyperlinks.Item(2).TextToDisplay = "9fb53365-d646-45a2-8ef4-49cdaa28c62b.md"
$zh.Hyperlinks.Item(3).TextToDisplay = "5fe23c8a-ea68-4f02-8d8b-6baca2a7de74.md"
$zh.Hyperlinks.Item(4).TextToDisplay = "5fe23c8a-ea68-4f02-8d8b-6baca2a7de74.md"

# Error Detail column widened to fit the long message now populated in P3.
$zh.Range("P1").EntireColumn.ColumnWidth = 39.17

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "9fb53365-d646-45a2-8ef4-49cdaa28c62b.md"
$de.Range("G2").Value = "9fb53365-d646-45a2-8ef4-49cdaa28c62b.fca86b0be0bc02ca2aa80548d8c5776cca403c63.de-de.xlf"
$de.Range("I2").Value = "9fb53365-d646-45a2-8ef4-49cdaa28c62b.md"
$de.Range("J2").Value = "9fb53365-d646-45a2-8ef4-49cdaa28c62b.fca86b0be0bc02ca2aa80548d8c5776cca403c63.de-de.xlf"

$de.Range("A3").Value = "5fe23c8a-ea68-4f02-8d8b-6baca2a7de74.md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("G3").Value = "5fe23c8a-ea68-4f02-8d8b-6baca2a7de74.e0e312a260014abee75e5e6bf1ab9bc9e18a2d7d.de-de.xlf"
$de.Range("H3").Value = "2016-09-03 06:51:40"
$de.Range("I3").Value = "5fe23c8a-ea68-4f02-8d8b-6baca2a7de74.md"
$de.Range("J3").Value = "5fe23c8a-ea68-4f02-8d8b-6baca2a7de74.e0e312a260014abee75e5e6bf1ab9bc9e18a2d7d.de-de.xlf"
$de.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/470930ec7e4338bec95c5b6e9ad25607e2184c5a/e2e/5fe23c8a-ea68-4f02-8d8b-6baca2a7de74.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7442d1f7ad367e65fbeed7c5c61faebaeb5310d8/e2e/5fe23c8a-ea68-4f02-8d8b-6baca2a7de74.md."

$de.Hyperlinks.Item(1).TextToDisplay = "9fb53365-d646-45a2-8ef4-49cdaa28c62b.md"
$de.Hyperlinks.Item(2).TextToDisplay = "9fb53365-d646-45a2-8ef4-49cdaa28c62b.md"
$de.Hyperlinks.Item(3).TextToDisplay = "5fe23c8a-ea68-4f02-8d8b-6baca2a7de74.md"
$de.Hyperlinks.Item(4).TextToDisplay = "5fe23c8a-ea68-4f02-8d8b-6baca2a7de74.md"

# Error Detail column widened to fit the long message now populated in P3.
$de.Range("P1").EntireColumn.ColumnWidth = 39.17
